# Rename the embedded Pearson/BTec logo pictures in every header & footer.
#
# The source document carries two distinct logos:
#   - "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
#     (embedded media/image2.png) appearing in both footers, originally
#     named "image2.png".
#   - "BTec_Logo-Orange" (embedded media/image1.jpg) appearing in both
#     headers, originally named "image1.jpg".
#
# The edit swaps each picture's display Name: the Pearson logo becomes
# "image1.png" and the BTec logo becomes "image2.jpg" (the wp:docPr / the
# InlineShape's Name, not the AlternativeText/descr, which is untouched).

$d = $word.ActiveDocument

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($h = 1; $h -le $sec.Headers.Count; $h++) {
        $hdr = $sec.Headers.Item($h)
        if (-not $hdr.Exists) { continue }
        for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }

    for ($f = 1; $f -le $sec.Footers.Count; $f++) {
        $ftr = $sec.Footers.Item($f)
        if (-not $ftr.Exists) { continue }
        for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
